$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# The subtitle placeholder currently holds two paragraphs:
#   1) "By Gavin Gunawardena"
#   2) "Code: <hyperlink>"
# Insert a new paragraph with the presentation date right after the
# first paragraph (and before the "Code:" paragraph), inheriting the
# same run formatting (Tahoma, 24pt) already used on that paragraph.
$byline = $tr.Paragraphs(1, 1)
$dateRange = $byline.InsertAfter([char]13 + "2/20/2023")
